# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the consolidated "全部类型" sheet to match the refreshed data pull.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 14949
$wsExhibit.Range("F3").Value = 18791
$wsExhibit.Range("F14").Value = 129
$wsExhibit.Range("F28").Value = 5997
$wsExhibit.Range("F30").Value = 69
$wsExhibit.Range("F31").Value = 165
$wsExhibit.Range("F34").Value = 5369
$wsExhibit.Range("F36").Value = 43

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 14949
$wsAll.Range("F3").Value = 18791
$wsAll.Range("F14").Value = 129
$wsAll.Range("F31").Value = 5997
$wsAll.Range("F33").Value = 69
$wsAll.Range("F34").Value = 165
$wsAll.Range("F37").Value = 5369
$wsAll.Range("F39").Value = 43
